$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$check = [char]0x2705

$oldLine1 = "$check 1000 Bs = 3.36 = 12733.02 pesos"
$newLine1 = "$check 1000 Bs = 3.34 = 12668.24 pesos"
$oldLine2 = "$check 12733.02 pesos = 3.33 = 966.49 Bs"
$newLine2 = "$check 12668.24 pesos = 3.32 = 975.78 Bs"

$text = $wsHoja1.Range("A1").Value2
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$wsHoja1.Range("A1").Value2 = $text

# --- Sheet "tasas": update the N10/O10/N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 299.497
$wsTasas.Range("O10").Value = 3794.1
$wsTasas.Range("N12").Value = 3814.3
$wsTasas.Range("O12").Value = 293.8
